$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 7 ("Looking at top rows"): new_df.head() -> df.head()
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(2)
$para7 = $shp7.TextFrame.TextRange.Paragraphs(1)
$para7.Text = "df.head()"
$para7.Characters(1, 7).Text = "df.head"
$para7.Characters(8, 2).Text = "()"

# the text box uses spAutoFit + vertical centering, so the box height/top
# shift slightly to keep the same vertical center once the text re-wraps
$shp7.Top = 379.04614
$shp7.Height = 155.90782

# ---------------------------------------------------------------------
# Slide 12 ("Series"): Numpy array with fancy index
# ---------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$shp12 = $s12.Shapes.Item(2)
$para12 = $shp12.TextFrame.TextRange.Paragraphs(2)
$para12.Text = "Numpy array with fancy index"
$para12.Characters(1, 5).Text = "Numpy"
$para12.Characters(6, 23).Text = " array with fancy index"

# ---------------------------------------------------------------------
# Slide 13 ("Boolean indexing"): df[df["SCORE"] > 30]
# ---------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$shp13 = $s13.Shapes.Item(2)
$para13 = $shp13.TextFrame.TextRange.Paragraphs(1)
$para13.Text = "df[df[“SCORE”] > 30]"
$para13.Characters(1, 2).Text = "df"
$para13.Characters(3, 1).Text = "["
$para13.Characters(4, 2).Text = "df"
$para13.Characters(6, 15).Text = "[“SCORE”] > 30]"

# re-splitting the run triggers the box's spAutoFit re-layout in this
# runtime; restore the unchanged box geometry from the source deck
$shp13.Top = 380.0
$shp13.Height = 154.00009

# ---------------------------------------------------------------------
# Slide 18 ("Concepts covered so Far"): drop the separate "dtype "
# bullet, turn the following bullet into ".head()"
# ---------------------------------------------------------------------
$s18 = $p.Slides.Item(18)
$shp18 = $s18.Shapes.Item(2)
$shp18.TextFrame.TextRange.Text = "Pandas`rDataFrames`r.head()`rBoolean indexing"
